# Update the simulation-output numeric columns (C:H) for rows 2-10 on the
# active sheet with refreshed values from a re-run of the model
# ("jan6 stable, with weird low kappa results").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, C, D, E, F, G, H
$data = @(
    @(2, 207.6771426899587, 670.7147553898792, 577.812872910898, 58.27642440795898, 0.2681797484556834, -12.27697058012694),
    @(3, 237.4366042710649, 777.0465877804547, 632.593492344443, 56.47904968261719, 0.3011502499381701, -8.791913332040506),
    @(4, 260.9055146073916, 769.0957211901207, 635.6250120368361, 55.71476745605469, 0.2891035795211792, -7.110600690419609),
    @(5, 272.3699359732539, 808.2246944543577, 649.313825610144, 55.15425109863281, 0.3258351335922877, -6.424721521794428),
    @(6, 277.627585341904, 811.7768842965512, 656.3094173192217, 54.79643249511719, 0.3240887582302093, -6.13164976385664),
    @(7, 385.4156814203054, 760.6041108193398, 666.8098509206569, 52.40727615356445, 0.2526679933071136, -1.328868571424224),
    @(8, 383.5409078367964, 758.7175875679043, 664.302173300609, 52.44415664672852, 0.2526491582393646, -1.308050690732231),
    @(9, 406.3610972076485, 778.1885662067625, 696.3223353070689, 50.11466979980469, 0.2197113633155823, -0.9550680434734532),
    @(10, 424.2297795461056, 840.2539961103175, 766.1450757336813, 46.74080276489258, 0.2077515870332718, -0.8646137478218536)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]   # C - cvar_05
    $ws.Cells.Item($row, 4).Value = $entry[2]   # D - expected_wealth
    $ws.Cells.Item($row, 5).Value = $entry[3]   # E - median_wealth
    $ws.Cells.Item($row, 6).Value = $entry[4]   # F - qsum_avg
    $ws.Cells.Item($row, 7).Value = $entry[5]   # G - p_med_avg
    $ws.Cells.Item($row, 8).Value = $entry[6]   # H - expected_LS
}
